# nowDate -> spotDate & whatDay 미사용
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 (planner_tbl): "name"/"작성자" column -> "userId"/"작성자 ID 값" column ---
# was: varchar(20), not null  ->  now: bigint, not null (length cleared)
$ws.Range("F9").Value = "작성자 ID 값"
$ws.Range("A9").Value = "userId"
$ws.Range("B9").Value = "bigint"
$ws.Range("C9").ClearContents()

# --- Row 21 (planner_spot_tbl): "nowDate" -> "spotDate" ---
$ws.Range("A21").Value = "spotDate"

# --- Row 22: whatDay column now marked as unused, note added in column G ---
$ws.Range("G22").Value = "일단 사용 x"
$ws.Range("G22").Font.Size = 18

# --- Selection moved from D15 to M22 ---
$ws.Range("M22").Select()

Write-Host "edits applied"
